$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7 (Ano 2025) with refreshed faturamento figures
$ws.Range("B7").Value = 2679493.91
$ws.Range("C7").Value = -39.69284721917066
$ws.Range("D7").Value = 2701
$ws.Range("E7").Value = 2701
$ws.Range("F7").Value = 992.0377304701963
$ws.Range("G7").Value = 5.744048711591176
